$wb = $excel.ActiveWorkbook

# Sheet1 = "图片路径" (Picture path) -- update the 5 picture paths to the
# new test_jianyuluntan upload paths (order 1,2,3,4,5) while leaving the
# header in A1 untouched.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "test_jianyuluntan/picture/1.jpg"
$ws1.Range("A3").Value = "test_jianyuluntan/picture/2.jpeg"
$ws1.Range("A4").Value = "test_jianyuluntan/picture/3.png"
$ws1.Range("A5").Value = "test_jianyuluntan/picture/4.gif"
$ws1.Range("A6").Value = "test_jianyuluntan/picture/5.jfif"

# Widen column A on sheet1 to fit the longer paths.
$ws1.Columns.Item(1).ColumnWidth = 34

# Move the selection on sheet1 to A5 and make it the active sheet/tab.
[void]$ws1.Range("A5").Select()
[void]$ws1.Activate()
